$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Project / PCB Version" header block (rows 1-4), which
# shifts the real table (header + data) up so it starts at row 1.
$ws.Rows("1:4").Delete()

# The diode's "Value" cell picked up a fuller description.
$ws.Range("B2").Value = "Diode Schottky 20V 0.5A"

# Column B grew a bit wider to fit the new text.
$ws.Columns("B").ColumnWidth = 21.109375

# Selection moved to B6 before saving.
$ws.Range("B6").Select()
